# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to "_FV2410" / "_FV2504"
# - Turn the data range into a proper Excel Table ("Table1")
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J -> "<name>_FV2410"
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = "$($oldHeaders[$i])_FV2410"
}

# Column K stays "diff"
$ws.Cells.Item(1, 11).Value2 = "diff"

# Columns L..U -> "<name>_FV2504"
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = "$($oldHeaders[$i])_FV2504"
}

# --- 2. Convert the A1:U88 range into an Excel Table -------------------
$tableRange = $ws.Range("A1:U88")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
